$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price strings so Excel does not
# auto-convert them to actual numbers (the source data stores every cell as text).
$textCells = @(
    "D4", "D5", "D6", "D8", "D9", "D11", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D26", "D27", "D30", "D31", "D32", "D35", "D36", "D37", "D38", "D41", "D43", "D44", "D45", "D47", "D48", "D49", "D50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.906.12"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").Value = "1.875.07"
$ws.Range("E3").Value = "  -0.93%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "0.7382"
$ws.Range("E5").Value = "  -4.66%  "

$ws.Range("D6").Value = "242.54"
$ws.Range("E6").Value = "  -0.55%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "0.3157"
$ws.Range("E8").Value = "  +0.85%  "

$ws.Range("D9").Value = "0.07184"
$ws.Range("E9").Value = "  -1.04%  "

$ws.Range("E10").Value = "  -4.35%  "

$ws.Range("D11").Value = "0.08417"
$ws.Range("E11").Value = "  -3.19%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.919.22"
$ws.Range("E12").Value = "  -3.46%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7507"
$ws.Range("E13").Value = "  -2.88%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.425"
$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").Value = "92.56"
$ws.Range("E15").Value = "  -2.08%  "

$ws.Range("D16").Value = "29.911.63"
$ws.Range("E16").Value = "  -0.58%  "

$ws.Range("D17").Value = "6.097"
$ws.Range("E17").Value = "  -1.89%  "

$ws.Range("D18").Value = "13.59"
$ws.Range("E18").Value = "  -2.59%  "

$ws.Range("D19").Value = "242.93"
$ws.Range("E19").Value = "  -1.15%  "

$ws.Range("D20").Value = "0.000007812"

$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").Value = "2.118.01"
$ws.Range("E22").Value = "  -6.48%  "

$ws.Range("E23").Value = "  -2.33%  "

$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("E25").Value = "  -3.13%  "

$ws.Range("D26").Value = "9.287"
$ws.Range("E26").Value = "  -2.70%  "

$ws.Range("D27").Value = "164.85"
$ws.Range("E27").Value = "  +1.13%  "

$ws.Range("E28").Value = "  -1.35%  "

$ws.Range("E29").Value = "  -0.69%  "

$ws.Range("D30").Value = "1.492"
$ws.Range("E30").Value = "  +4.23%  "

$ws.Range("D31").Value = "4.603"
$ws.Range("E31").Value = "  +1.54%  "

$ws.Range("D32").Value = "1.530"
$ws.Range("E32").Value = "  -1.03%  "

$ws.Range("E33").Value = "  +2.83%  "

$ws.Range("E34").Value = "  -2.56%  "

$ws.Range("D35").Value = "1.235"
$ws.Range("E35").Value = "  -1.16%  "

$ws.Range("D36").Value = "0.7542"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").Value = "0.9999"
$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D38").Value = "2.696"
$ws.Range("E38").Value = "  +0.23%  "

$ws.Range("E39").Value = "  -1.31%  "

$ws.Range("E40").Value = "  -1.07%  "

$ws.Range("D41").Value = "0.4518"
$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").Value = "1.110.60"
$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("D43").Value = "6.060"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").Value = "72.37"
$ws.Range("E44").Value = "  -1.95%  "

$ws.Range("D45").Value = "0.8573"
$ws.Range("E45").Value = "  +0.40%  "

$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("D47").Value = "103.41"
$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").Value = "7.655"
$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("D49").Value = "3.089"
$ws.Range("E49").Value = "  +2.86%  "

$ws.Range("D50").Value = "1.839"
$ws.Range("E50").Value = "  -2.57%  "

$ws.Range("D51").Value = "2.016.25"
$ws.Range("E51").Value = "  -7.99%  "
